# Update the cryptos list with refreshed prices / volume percentages.
# Note: several "price" strings (column D) look like plain numbers to Excel
# (e.g. "94.11"), which would otherwise be auto-converted to a numeric value
# (losing formatting / trailing zeros, e.g. "1.90" -> 1.9). A leading
# apostrophe forces those assignments to be stored as text, matching the
# original inline-string cell type, exactly like typing them in the Excel UI.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.215.47'
$ws.Range('E2').Value = '  -3.57%  '
$ws.Range('D3').Value = '2.465.47'
$ws.Range('E3').Value = '  -2.47%  '
$ws.Range('E5').Value = '  +0.71%  '
$ws.Range('D6').Value = '''94.11'
$ws.Range('E6').Value = '  -6.37%  '
$ws.Range('E7').Value = '  -2.87%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').Value = '''0.497'
$ws.Range('E9').Value = '  -4.57%  '
$ws.Range('D10').Value = '''33.37'
$ws.Range('E10').Value = '  -6.23%  '
$ws.Range('E11').Value = '  -3.32%  '
$ws.Range('E12').Value = '  -1.25%  '
$ws.Range('D13').Value = '''6.99'
$ws.Range('E13').Value = '  -4.39%  '
$ws.Range('D14').Value = '2.844.95'
$ws.Range('E14').Value = '  -2.48%  '
$ws.Range('D15').Value = '2.475.19'
$ws.Range('E15').Value = '  -3.11%  '
$ws.Range('D16').Value = '''14.81'
$ws.Range('E16').Value = '  -3.29%  '
$ws.Range('D17').Value = '''0.783'
$ws.Range('E17').Value = '  -3.58%  '
$ws.Range('D18').Value = '41.178.02'
$ws.Range('E18').Value = '  -3.62%  '
$ws.Range('E19').Value = '  -5.64%  '
$ws.Range('E20').Value = '  -2.99%  '
$ws.Range('D21').Value = '''11.29'
$ws.Range('E21').Value = '  -7.91%  '
$ws.Range('D22').Value = '''68.34'
$ws.Range('E22').Value = '  -1.41%  '
$ws.Range('D23').Value = '''235.62'
$ws.Range('E23').Value = '  -2.99%  '
$ws.Range('E24').Value = '  -3.98%  '
$ws.Range('D26').Value = '''1.90'
$ws.Range('E26').Value = '  -5.90%  '
$ws.Range('D27').Value = '''23.96'
$ws.Range('E27').Value = '  -5.69%  '
$ws.Range('E28').Value = '  -6.27%  '
$ws.Range('E29').Value = '  -5.35%  '
$ws.Range('D30').Value = '''36.36'
$ws.Range('E30').Value = '  -5.39%  '
$ws.Range('D31').Value = '''152.79'
$ws.Range('E31').Value = '  -4.41%  '
$ws.Range('D32').Value = '''5.48'
$ws.Range('E32').Value = '  -4.96%  '
$ws.Range('E33').Value = '  -4.96%  '
$ws.Range('E34').Value = '  -3.49%  '
$ws.Range('D35').Value = '''0.0744'
$ws.Range('E35').Value = '  -4.98%  '
$ws.Range('D36').Value = '''3.04'
$ws.Range('E36').Value = '  -2.45%  '
$ws.Range('B37').Value = 'ARBITRUM'
$ws.Range('C37').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D37').Value = '''1.88'
$ws.Range('E37').Value = '  -4.14%  '
$ws.Range('B38').Value = 'Celestia'
$ws.Range('C38').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D38').Value = '''17.01'
$ws.Range('E38').Value = '  -7.99%  '
$ws.Range('E39').Value = '  -2.95%  '
$ws.Range('E40').Value = '  -7.46%  '
$ws.Range('D41').Value = '''4.19'
$ws.Range('E41').Value = '  -0.86%  '
$ws.Range('E42').Value = '  +0.13%  '
$ws.Range('D43').Value = '''20.18'
$ws.Range('E43').Value = '  -9.45%  '
$ws.Range('D44').Value = '1.967.16'
$ws.Range('E44').Value = '  -1.40%  '
$ws.Range('E45').Value = '  -5.38%  '
$ws.Range('D46').Value = '''3.04'
$ws.Range('E46').Value = '  -7.99%  '
$ws.Range('E47').Value = '  -2.08%  '
$ws.Range('D48').Value = '''69.37'
$ws.Range('E48').Value = '  -3.09%  '
$ws.Range('D49').Value = '''97.19'
$ws.Range('E49').Value = '  -3.61%  '
$ws.Range('E50').Value = '  -6.44%  '
$ws.Range('E51').Value = '  -6.84%  '
